$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 987.80646
$ws.Range("I33").Value = 848.8261
$ws.Range("J33").Value = 1387.375
$ws.Range("K33").Value = 848.8261
$ws.Range("L33").Value = 1387.375
$ws.Range("M33").Value = -619.8261
$ws.Range("N33").Value = -1845.375

$ws.Range("H123").Value = 30000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 30000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800

$ws.Range("H132").Value = 3197.3489
$ws.Range("I132").Value = 3035.4048
$ws.Range("J132").Value = 9999
$ws.Range("K132").Value = 9106.214399999999
$ws.Range("L132").Value = 29997
$ws.Range("M132").Value = -6576.214399999999
$ws.Range("N132").Value = -35057

$ws.Range("H135").Value = 3872.2856
$ws.Range("I135").Value = 3351
$ws.Range("J135").Value = 7000
$ws.Range("K135").Value = 30159
$ws.Range("L135").Value = 63000
$ws.Range("M135").Value = -27624
$ws.Range("N135").Value = -68070

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()

$ws.Range("H32").Value = 420067.34
$ws.Range("I32").Value = 493171.38
$ws.Range("J32").Value = 17995.166
$ws.Range("K32").Value = 493171.38
$ws.Range("L32").Value = 17995.166
$ws.Range("M32").Value = -492884.38
$ws.Range("N32").Value = -18569.166

$ws.Range("H45").Value = 2625.7778
$ws.Range("I45").Value = 2019.4286
$ws.Range("J45").Value = 3278.7693
$ws.Range("K45").Value = 2019.4286
$ws.Range("L45").Value = 3278.7693
$ws.Range("M45").Value = -1642.4286
$ws.Range("N45").Value = -4032.7693

$ws.Range("H92").Value = 50895.715
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 50895.715
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 50895.715
$ws.Range("N92").Value = -55887.715

$ws.Range("H132").Value = 3604.9697
$ws.Range("I132").Value = 3477.75
$ws.Range("J132").Value = 3944.2222
$ws.Range("K132").Value = 10433.25
$ws.Range("L132").Value = 11832.6666
$ws.Range("M132").Value = -7903.25
$ws.Range("N132").Value = -16892.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 65000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 65000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 65000
$ws.Range("N51").Value = -65982

$ws.Range("H134").Value = 2685.1458
$ws.Range("I134").Value = 2612.4473
$ws.Range("J134").Value = 2961.4
$ws.Range("K134").Value = 7837.341899999999
$ws.Range("L134").Value = 8884.200000000001
$ws.Range("M134").Value = -5302.341899999999
$ws.Range("N134").Value = -13954.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4515.5947
$ws.Range("I31").Value = 1094.0435
$ws.Range("J31").Value = 10136.714
$ws.Range("K31").Value = 1094.0435
$ws.Range("L31").Value = 10136.714
$ws.Range("M31").Value = -799.0435
$ws.Range("N31").Value = -10726.714

$ws.Range("H34").Value = 4515.5947
$ws.Range("I34").Value = 1094.0435
$ws.Range("J34").Value = 10136.714
$ws.Range("K34").Value = 1094.0435
$ws.Range("L34").Value = 10136.714
$ws.Range("M34").Value = -892.0435
$ws.Range("N34").Value = -10540.714

$ws.Range("H99").Value = 1940.0435
$ws.Range("I99").Value = 1804.2
$ws.Range("J99").Value = 1977.7778
$ws.Range("K99").Value = 1804.2
$ws.Range("L99").Value = 1977.7778
$ws.Range("M99").Value = -306.2
$ws.Range("N99").Value = -4973.7778

$ws.Range("H100").Value = 49663.332
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 49663.332
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 49663.332
$ws.Range("N100").Value = -51827.332

$ws.Range("H126").Value = 1940.0435
$ws.Range("I126").Value = 1804.2
$ws.Range("J126").Value = 1977.7778
$ws.Range("K126").Value = 5412.6
$ws.Range("L126").Value = 5933.3334
$ws.Range("M126").Value = -2942.6
$ws.Range("N126").Value = -10873.3334

$ws.Range("H132").Value = 4763305
$ws.Range("I132").Value = 855.5925999999999
$ws.Range("J132").Value = 20836572
$ws.Range("K132").Value = 2566.7778
$ws.Range("L132").Value = 62509716
$ws.Range("M132").Value = -36.77779999999984
$ws.Range("N132").Value = -62514776

$ws.Range("H134").Value = 1779
$ws.Range("I134").Value = 1627.25
$ws.Range("J134").Value = 2386
$ws.Range("K134").Value = 4881.75
$ws.Range("L134").Value = 7158
$ws.Range("M134").Value = -2346.75
$ws.Range("N134").Value = -12228

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8752637
$ws.Range("I4").Value = 8002618
$ws.Range("J4").Value = 9093555
$ws.Range("K4").Value = 24007854
$ws.Range("L4").Value = 27280665
$ws.Range("M4").Value = -24007742
$ws.Range("N4").Value = -27280889

$ws.Range("H39").Value = 1279.3529
$ws.Range("I39").Value = 447.83334
$ws.Range("J39").Value = 1732.909
$ws.Range("K39").Value = 1343.50002
$ws.Range("L39").Value = 5198.727000000001
$ws.Range("M39").Value = -1049.50002
$ws.Range("N39").Value = -5786.727000000001

$ws.Range("H113").Value = 688.0909
$ws.Range("I113").Value = 409
$ws.Range("J113").Value = 1286.1428
$ws.Range("K113").Value = 1227
$ws.Range("L113").Value = 3858.4284
$ws.Range("M113").Value = 943
$ws.Range("N113").Value = -8198.428400000001

$ws.Range("H122").Value = 7758.0713
$ws.Range("I122").Value = 551.25
$ws.Range("J122").Value = 50999
$ws.Range("K122").Value = 4961.25
$ws.Range("L122").Value = 458991
$ws.Range("M122").Value = -2511.25
$ws.Range("N122").Value = -463891

$ws.Range("H137").Value = 14610.223
$ws.Range("I137").Value = 19776.5
$ws.Range("J137").Value = 4277.6665
$ws.Range("K137").Value = 59329.5
$ws.Range("L137").Value = 12832.9995
$ws.Range("M137").Value = -54229.5
$ws.Range("N137").Value = -23032.9995

$ws.Range("H140").Value = 1422.8462
$ws.Range("I140").Value = 956.6957
$ws.Range("J140").Value = 4996.6665
$ws.Range("K140").Value = 2870.0871
$ws.Range("L140").Value = 14989.9995
$ws.Range("M140").Value = 2309.9129
$ws.Range("N140").Value = -25349.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1910.0555
$ws.Range("I97").Value = 1540
$ws.Range("J97").Value = 2650.1667
$ws.Range("K97").Value = 1540
$ws.Range("L97").Value = 2650.1667
$ws.Range("M97").Value = -1044
$ws.Range("N97").Value = -3642.1667

$ws.Range("H126").Value = 3529.7273
$ws.Range("I126").Value = 2982.8
$ws.Range("J126").Value = 3985.5
$ws.Range("K126").Value = 8948.400000000001
$ws.Range("L126").Value = 11956.5
$ws.Range("M126").Value = -6478.400000000001
$ws.Range("N126").Value = -16896.5

$ws.Range("H132").Value = 2111.6956
$ws.Range("I132").Value = 1795.7894
$ws.Range("J132").Value = 3612.25
$ws.Range("K132").Value = 5387.3682
$ws.Range("L132").Value = 10836.75
$ws.Range("M132").Value = -2857.3682
$ws.Range("N132").Value = -15896.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4603.684
$ws.Range("I40").Value = 1530
$ws.Range("J40").Value = 5180
$ws.Range("K40").Value = 1530
$ws.Range("L40").Value = 5180
$ws.Range("M40").Value = -1394
$ws.Range("N40").Value = -5452

$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()

$ws.Range("H122").Value = 2630
$ws.Range("I122").Value = 2545
$ws.Range("J122").Value = 2800
$ws.Range("K122").Value = 7635
$ws.Range("L122").Value = 8400
$ws.Range("M122").Value = -5185
$ws.Range("N122").Value = -13300

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 77400
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 77400
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 77400
$ws.Range("N62").Value = -78648

$ws.Range("H65").Value = 77400
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 77400
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 387000
$ws.Range("N65").Value = -393240

$ws.Range("H107").Value = 743.625
$ws.Range("I107").Value = 717.25
$ws.Range("J107").Value = 770
$ws.Range("K107").Value = 2151.75
$ws.Range("L107").Value = 2310
$ws.Range("M107").Value = -231.75
$ws.Range("N107").Value = -6150

$ws.Range("H122").Value = 3525.2
$ws.Range("I122").Value = 1316.6666
$ws.Range("J122").Value = 4471.7144
$ws.Range("K122").Value = 3949.9998
$ws.Range("L122").Value = 13415.1432
$ws.Range("M122").Value = -1499.9998
$ws.Range("N122").Value = -18315.1432

$ws.Range("H126").Value = 1596.1333
$ws.Range("I126").Value = 1395.2
$ws.Range("J126").Value = 1998
$ws.Range("K126").Value = 4185.6
$ws.Range("L126").Value = 5994
$ws.Range("M126").Value = -1715.6
$ws.Range("N126").Value = -10934

$ws.Range("H132").Value = 5053254.5
$ws.Range("I132").Value = 4011.1667
$ws.Range("J132").Value = 7938536.5
$ws.Range("K132").Value = 12033.5001
$ws.Range("L132").Value = 23815609.5
$ws.Range("M132").Value = -9503.500100000001
$ws.Range("N132").Value = -23820669.5

$ws.Range("H136").Value = 2309.1956
$ws.Range("I136").Value = 2024.3928
$ws.Range("J136").Value = 2752.2222
$ws.Range("K136").Value = 6073.178400000001
$ws.Range("L136").Value = 8256.6666
$ws.Range("M136").Value = -3523.178400000001
$ws.Range("N136").Value = -13356.6666
